$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 5.005600929260254
$ws.Range("B1").Value = 3.635507106781006
$ws.Range("C1").Value = 2.743049144744873
$ws.Range("D1").Value = 2.304853677749634
$ws.Range("E1").Value = 2.2149817943573
